$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 14248.083
$ws.Range("I51").Value = 41749.5
$ws.Range("K51").Value = 41749.5
$ws.Range("M51").Value = -41265.5

$ws.Range("H86").Value = 44543544
$ws.Range("I86").Value = 53746508
$ws.Range("K86").Value = 53746508
$ws.Range("M86").Value = -53745385

$ws.Range("H89").Value = 44543544
$ws.Range("I89").Value = 53746508
$ws.Range("K89").Value = 268732540
$ws.Range("M89").Value = -268726924

$ws.Range("H116").Value = 62510750
$ws.Range("I116").Value = 250000000
$ws.Range("J116").Value = 14333.333
$ws.Range("K116").Value = 250000000
$ws.Range("L116").Value = 14333.333
$ws.Range("M116").Value = -249996558
$ws.Range("N116").Value = -21217.333

$ws.Range("H138").Value = 5317.6226
$ws.Range("J138").Value = 7898.1875
$ws.Range("L138").Value = 23694.5625
$ws.Range("N138").Value = -33974.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2240652
$ws.Range("I32").Value = 2457980.5
$ws.Range("J32").Value = 23902.4
$ws.Range("K32").Value = 2457980.5
$ws.Range("L32").Value = 23902.4
$ws.Range("M32").Value = -2457693.5
$ws.Range("N32").Value = -24476.4

$ws.Range("H63").Value = 1051.7778
$ws.Range("I63").Value = 1017.75
$ws.Range("J63").Value = 1079
$ws.Range("K63").Value = 1017.75
$ws.Range("L63").Value = 1079
$ws.Range("M63").Value = -331.75
$ws.Range("N63").Value = -2451

$ws.Range("H66").Value = 1051.7778
$ws.Range("I66").Value = 1017.75
$ws.Range("J66").Value = 1079
$ws.Range("K66").Value = 5088.75
$ws.Range("L66").Value = 5395
$ws.Range("M66").Value = -1656.75
$ws.Range("N66").Value = -12259

$ws.Range("H97").Value = 4630611
$ws.Range("I97").Value = 910.8182
$ws.Range("K97").Value = 910.8182
$ws.Range("M97").Value = -414.8182

$ws.Range("H122").Value = 3107.76
$ws.Range("I122").Value = 1795
$ws.Range("K122").Value = 5385
$ws.Range("M122").Value = -2935

$ws.Range("H132").Value = 4608.5396
$ws.Range("I132").Value = 3350.75
$ws.Range("J132").Value = 6796
$ws.Range("K132").Value = 10052.25
$ws.Range("L132").Value = 20388
$ws.Range("M132").Value = -7522.25
$ws.Range("N132").Value = -25448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12822890
$ws.Range("I20").Value = 18520924
$ws.Range("J20").Value = 2313.75
$ws.Range("K20").Value = 18520924
$ws.Range("L20").Value = 2313.75
$ws.Range("M20").Value = -18520677
$ws.Range("N20").Value = -2807.75

$ws.Range("H94").Value = 1611.3096
$ws.Range("J94").Value = 2939.2354
$ws.Range("L94").Value = 2939.2354
$ws.Range("N94").Value = -3841.2354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6235.0312
$ws.Range("I31").Value = 2759.2104
$ws.Range("K31").Value = 2759.2104
$ws.Range("M31").Value = -2464.2104

$ws.Range("H34").Value = 6235.0312
$ws.Range("I34").Value = 2759.2104
$ws.Range("K34").Value = 2759.2104
$ws.Range("M34").Value = -2557.2104

$ws.Range("H58").Value = 14292977
$ws.Range("I58").Value = 31252520
$ws.Range("K58").Value = 31252520
$ws.Range("M58").Value = -31252317

$ws.Range("H132").Value = 5085.0786
$ws.Range("I132").Value = 2482.0881
$ws.Range("J132").Value = 10291.059
$ws.Range("K132").Value = 7446.2643
$ws.Range("L132").Value = 30873.177
$ws.Range("M132").Value = -4916.2643
$ws.Range("N132").Value = -35933.177

$ws.Range("H136").Value = 14292977
$ws.Range("I136").Value = 31252520
$ws.Range("K136").Value = 93757560
$ws.Range("M136").Value = -93755010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 122478210
$ws.Range("I9").Value = 66681970
$ws.Range("K9").Value = 200045910
$ws.Range("M9").Value = -200045686

$ws.Range("H18").Value = 187.5
$ws.Range("I18").Value = 197.55556
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = 592.66668
$ws.Range("L18").Value = 291
$ws.Range("M18").Value = -423.66668
$ws.Range("N18").Value = -629

$ws.Range("H131").Value = 1650.2858
$ws.Range("I131").Value = 1395.909
$ws.Range("J131").Value = 1930.1
$ws.Range("K131").Value = 4187.727000000001
$ws.Range("L131").Value = 5790.299999999999
$ws.Range("M131").Value = 852.2729999999992
$ws.Range("N131").Value = -15870.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 810001.7
$ws.Range("I70").Value = 8000000
$ws.Range("J70").Value = 11113
$ws.Range("K70").Value = 8000000
$ws.Range("L70").Value = 11113
$ws.Range("M70").Value = -7999730
$ws.Range("N70").Value = -11653

$ws.Range("H73").Value = 810001.7
$ws.Range("I73").Value = 8000000
$ws.Range("J73").Value = 11113
$ws.Range("K73").Value = 8000000
$ws.Range("L73").Value = 11113
$ws.Range("M73").Value = -7999064
$ws.Range("N73").Value = -12985

$ws.Range("H105").Value = 30335.5
$ws.Range("J105").Value = 30335.5
$ws.Range("L105").Value = 30335.5
$ws.Range("N105").Value = -37323.5

$ws.Range("H107").Value = 1000164.75
$ws.Range("J107").Value = 77.5
$ws.Range("L107").Value = 77.5
$ws.Range("N107").Value = -3917.5

$ws.Range("H132").Value = 4647.3105
$ws.Range("J132").Value = 10094.056
$ws.Range("L132").Value = 30282.168
$ws.Range("N132").Value = -35342.16800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 899996.3
$ws.Range("I20").Value = 99995
$ws.Range("J20").Value = 1299997
$ws.Range("K20").Value = 99995
$ws.Range("L20").Value = 1299997
$ws.Range("M20").Value = -99769
$ws.Range("N20").Value = -1300449

$ws.Range("H22").Value = 2070.4546
$ws.Range("I22").Value = 990.3333
$ws.Range("J22").Value = 2475.5
$ws.Range("K22").Value = 990.3333
$ws.Range("L22").Value = 2475.5
$ws.Range("M22").Value = -695.3333
$ws.Range("N22").Value = -3065.5

$ws.Range("H27").Value = 2070.4546
$ws.Range("I27").Value = 990.3333
$ws.Range("J27").Value = 2475.5
$ws.Range("K27").Value = 990.3333
$ws.Range("L27").Value = 2475.5
$ws.Range("M27").Value = -883.3333
$ws.Range("N27").Value = -2689.5

$ws.Range("H46").Value = 5294896
$ws.Range("I46").Value = 2000.3334
$ws.Range("J46").Value = 6177045.5
$ws.Range("K46").Value = 2000.3334
$ws.Range("L46").Value = 6177045.5
$ws.Range("M46").Value = -1812.3334
$ws.Range("N46").Value = -6177421.5

$ws.Range("H61").Value = 8252.272000000001
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8252.272000000001
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8252.272000000001
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -8656.272000000001

$ws.Range("H68").Value = 4376.4287
$ws.Range("I68").Value = 2673.4546
$ws.Range("J68").Value = 6249.7
$ws.Range("K68").Value = 2673.4546
$ws.Range("L68").Value = 6249.7
$ws.Range("M68").Value = -1924.4546
$ws.Range("N68").Value = -7747.7

$ws.Range("H71").Value = 4376.4287
$ws.Range("I71").Value = 2673.4546
$ws.Range("J71").Value = 6249.7
$ws.Range("K71").Value = 13367.273
$ws.Range("L71").Value = 31248.5
$ws.Range("M71").Value = -9623.273000000001
$ws.Range("N71").Value = -38736.5

$ws.Range("H100").Value = 4645.6
$ws.Range("I100").Value = 4222
$ws.Range("J100").Value = 4928
$ws.Range("K100").Value = 4222
$ws.Range("L100").Value = 4928
$ws.Range("M100").Value = -3681
$ws.Range("N100").Value = -6010

$ws.Range("H113").Value = 8252.272000000001
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8252.272000000001
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8252.272000000001
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -12592.272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 34999.5
$ws.Range("I15").Value = 34999.5
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 34999.5
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -34711.5
$ws.Range("N15").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H62").Value = 50061.547
$ws.Range("I62").Value = 66972.875
$ws.Range("J62").Value = 4964.6665
$ws.Range("K62").Value = 66972.875
$ws.Range("L62").Value = 4964.6665
$ws.Range("M62").Value = -66348.875
$ws.Range("N62").Value = -6212.6665

$ws.Range("H65").Value = 50061.547
$ws.Range("I65").Value = 66972.875
$ws.Range("J65").Value = 4964.6665
$ws.Range("K65").Value = 334864.375
$ws.Range("L65").Value = 24823.3325
$ws.Range("M65").Value = -331744.375
$ws.Range("N65").Value = -31063.3325

$ws.Range("H81").Value = 22245488
$ws.Range("I81").Value = 1599.75
$ws.Range("J81").Value = 40040600
$ws.Range("K81").Value = 3199.5
$ws.Range("L81").Value = 80081200
$ws.Range("M81").Value = -2138.5
$ws.Range("N81").Value = -80083322

$ws.Range("H84").Value = 22245488
$ws.Range("I84").Value = 1599.75
$ws.Range("J84").Value = 40040600
$ws.Range("K84").Value = 15997.5
$ws.Range("L84").Value = 400406000
$ws.Range("M84").Value = -10693.5
$ws.Range("N84").Value = -400416608

$ws.Range("H96").Value = 1171.5555
$ws.Range("I96").Value = 1419.8
$ws.Range("J96").Value = 861.25
$ws.Range("K96").Value = 1419.8
$ws.Range("L96").Value = 861.25
$ws.Range("M96").Value = -46.79999999999995
$ws.Range("N96").Value = -3607.25

$ws.Range("H105").Value = 38331.668
$ws.Range("J105").Value = 38331.668
$ws.Range("L105").Value = 38331.668
$ws.Range("N105").Value = -45319.668

$ws.Range("H107").Value = 33334920
$ws.Range("I107").Value = 2564
$ws.Range("J107").Value = 47620216
$ws.Range("K107").Value = 7692
$ws.Range("L107").Value = 142860648
$ws.Range("M107").Value = -5772
$ws.Range("N107").Value = -142864488

$ws.Range("H113").Value = 12666.477
$ws.Range("I113").Value = 15553.059
$ws.Range("K113").Value = 46659.177
$ws.Range("M113").Value = -44489.177
